# Corrección del error 53
# Actualización de la historia 64 al eliminar botón regresar
#
# En la tabla de la historia de usuario "lista de ventas", tras el
# párrafo que describe el botón "Ver Detalles", existía un párrafo
# adicional que documentaba un botón "Regresar" que ya no aplica.
# Se elimina ese párrafo completo (texto + marca de párrafo), dejando
# intacto el resto de la celda y de la tabla.

$d = $word.ActiveDocument

$openQuote  = [char]0x201C
$closeQuote = [char]0x201D
$target = "Bot" + [char]0x00F3 + "n " + $openQuote + "Regresar" + $closeQuote + `
          ": Creaci" + [char]0x00F3 + "n de un bot" + [char]0x00F3 + `
          "n que nos redirija a la pantalla principal."

$deleted = $false

# Localizar el párrafo exacto (el texto es único en todo el documento) y
# eliminar su rango completo, lo que también retira la marca de párrafo
# y por tanto la entrada completa <w:p> del XML subyacente.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Contains($target)) {
        $p.Range.Delete()
        $deleted = $true
        break
    }
}

if (-not $deleted) {
    # Alternativa de respaldo usando Find, por si la iteración anterior no
    # localizara el párrafo.
    $rng = $d.Content
    $found = $rng.Find.Execute($target, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if ($found) {
        $rng.Paragraphs(1).Range.Delete()
    }
}
